$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.524.29"
$ws.Range("D3").Value = "2.468.00"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'314.54"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'91.94"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "'32.27"
$ws.Range("E10").Value = "  -4.27%  "
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D13").Value = "2.848.23"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "'6.84"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("E15").Value = "  +2.08%  "
$ws.Range("D16").Value = "2.445.43"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "'0.774"
$ws.Range("D18").Value = "41.513.49"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").Value = "'11.07"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "'235.96"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "'9.69"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "'35.38"
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("D31").Value = "'155.87"
$ws.Range("E31").Value = "  +2.13%  "
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "'0.0758"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").Value = "'17.31"
$ws.Range("E35").Value = "  -4.85%  "
$ws.Range("E36").Value = "  -7.48%  "
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -6.27%  "
$ws.Range("D40").Value = "'2.22"
$ws.Range("E40").Value = "  -10.82%  "
$ws.Range("E41").Value = "  -5.57%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").Value = "1.945.94"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").Value = "'18.42"
$ws.Range("E45").Value = "  -7.90%  "
$ws.Range("D46").Value = "'2.93"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "'9.03"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("D48").Value = "2.706.57"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "'96.85"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").Value = "'67.11"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").Value = "'52.18"
$ws.Range("E51").Value = "  +1.89%  "
